$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unit "assignment/quiz" -> "assignment" wording updates, plus final-project wording tweaks
$ws.Range("E7").Value = "b. Preview Unit 1 assignment"
$ws.Range("E8").Value = "a. Review Unit 1 assignment`n"
$ws.Range("E11").Value = "b. Preview Unit 2 assignment and final project"
$ws.Range("E12").Value = "a. Review Unit 2 assignment"
$ws.Range("E15").Value = "b. Preview Unit 3 assignment"
$ws.Range("E17").Value = "b. Prepare for final project"
$ws.Range("E18").Value = "a. Review Unit 3 assignment`n"
$ws.Range("E20").Value = "c. Preview Unit 4 assignment"
$ws.Range("E21").Value = "a. Review Unit 4 assignment"

# Update the active selection to match the saved view state
$ws.Range("E19").Select()
